$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header cells J1:O1, styled like the existing header row (copy from I1) ---
$headers = @("IsEdited", "IsDeleted", "AttachmentType", "AttachmentName", "AttachmentPath", "AttachmentSize")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 10 + $i   # J=10 .. O=15
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}
# Copy the header style (bold/border/fill/alignment) from I1 onto J1:O1
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:O1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# --- Existing data rows 2 & 3: new trailing columns are empty strings ---
foreach ($r in 2, 3) {
    for ($c = 10; $c -le 15; $c++) {
        $ws.Cells.Item($r, $c).Value = ""
    }
}

# --- New row 4 data ---
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "U001"
$ws.Cells.Item(4, 3).Value = "آقای گلستانی"
$ws.Cells.Item(4, 4).Value = "سلام"
$ws.Cells.Item(4, 5).Value = "2026-02-03 08:56:56"

# JalaliDate ("1404/11/14") looks like a date to Excel's smart entry, so
# temporarily force Text format to keep it verbatim, then clear the format
# back off again so the cell does not end up carrying a stray style index.
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "1404/11/14"
$ws.Cells.Item(4, 6).ClearFormats()

$ws.Cells.Item(4, 7).Value = "08:56"
$ws.Cells.Item(4, 8).Value = $false
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = $false
$ws.Cells.Item(4, 11).Value = $false
$ws.Cells.Item(4, 12).Value = "none"
$ws.Cells.Item(4, 13).Value = ""
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 15).Value = ""
